# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every value in column C (rows 2-396) moves from serial date 45202
# (2023-10-03) to serial date 45203 (2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 396 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
